# Feedback on first version of the production Definition
#
# On the "Hardware Development Process" sheet, a new deliverable row
# ("Interface list") is inserted right before the existing "Major
# Components BOM" row (old row 10), pushing every row below it down by
# one. The new row is populated the same way as its sibling deliverable
# rows (e.g. row 9, "Behavior Definition"): estimated hours, a due date,
# an assigned/completion date, and the "End of Day" note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hardware Development Process")

# Insert a new blank row at position 10 (existing row 10 and everything
# below shifts down by one row); formatting is inherited from the row
# above, matching the sibling deliverable rows.
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 2).Value = "Interface list"
$ws.Cells.Item(10, 3).Value = 1.5
$ws.Cells.Item(10, 4).Value = 42872
$ws.Cells.Item(10, 5).Value = 42874
$ws.Cells.Item(10, 6).Value = "End of Day"

# Reflect the author's final cursor position on the sheet.
$ws.Range("C17").Select()
